$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.026.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.88%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.830.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +0.04%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''1.008'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.32%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''311.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -0.78%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''1.008'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -0.19%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.4656'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -1.55%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.3715'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +0.97%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.07435'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.04%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.8725'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -1.33%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''20.05'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -2.04%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.07869'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = '''1.834.45'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -4.22%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''6.603'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +0.90%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''5.372'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -0.84%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''92.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -1.58%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.12%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.000008995'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +2.35%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = '''  -0.27%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''14.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -0.09%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''27.056.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -2.25%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '''  -2.12%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''10.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -0.22%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.064.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -2.51%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''152.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = '''1.838'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -3.26%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''18.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -1.86%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''2.101'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -1.50%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''5.131'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -1.83%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''115.86'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -0.98%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.08884'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -1.05%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''2.975'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +0.92%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = '''ImmutableX'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = '''0.7300'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -2.11%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = '''Filecoin'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = '''4.456'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -1.65%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''1.136'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -3.20%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''2.483'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +2.80%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''1.080'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -1.16%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = '''  +0.00%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = '''FraxShare'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = '''7.388'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +2.34%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = '''Hedera'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = '''0.05250'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -1.56%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''2.927'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.81%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.5180'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -1.84%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.1631'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -1.36%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''0.8581'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = '''8.245'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -2.72%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''0.4861'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -0.76%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = '''PaxDollar'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''1.008'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.19%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = '''EnergySwap'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''10.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -3.14%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''102.78'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -2.11%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''1.627'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -2.06%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.06251'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -0.77%  '
$ws.Range('E51').Style = 'Normal'
